# "remove duplicate future page"
#
# The deck has two adjacent "Future Expectations" slides (a leftover
# duplicate). Delete the first of the two (slide 8 - SlideID 264,
# "Free app for consumers..." body) and keep the second one (now
# shifted up to slide 8, SlideID 279, "Consumers / Retailer..." body).

$p = $ppt.ActivePresentation

# Find the duplicate "Future Expectations" slide that precedes the
# "Technical Challenges" slide - i.e. the first of the two consecutive
# "Future Expectations" slides - and remove it.
$targetIndex = -1
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    $title = $s.Shapes.Item(1).TextFrame.TextRange.Text
    if ($title -eq "Future Expectations") {
        $next = $p.Slides.Item($i + 1)
        $nextTitle = $next.Shapes.Item(1).TextFrame.TextRange.Text
        if ($nextTitle -eq "Future Expectations") {
            $targetIndex = $i
            break
        }
    }
}

if ($targetIndex -gt 0) {
    $p.Slides.Item($targetIndex).Delete()
}
